$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn off the binary switches B7:B15 (was all 1, now all 0) ---
$ws.Range("B7:B15").Value = 0

# --- Add new "j" letter -> 7-segment lookup block, rows 25:29 ---
$ws.Range("A25").Value = "j"

# --- Add new XNOR logic-gate block at G10:H11 ---
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H10").Value = "XNOR"
$ws.Range("H11").Formula = "=IF(G10=G11, 1, 0)"

$ws.Range("B25").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "0245689abcdefghjklmnopqrstuvwxyz" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("C25").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "02356789abcdefgjklmnopqrstuvwxyz" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("D25").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "0123456789abcefghijkmnopqrstuvwxyz" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"

$ws.Range("A26").Font.Underline = $true
$ws.Range("B26").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "045689abcdefghklmnopqrsuvwxy" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("D26").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "01234789abdh" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"

$ws.Range("A27").HorizontalAlignment = -4108
$ws.Range("B27").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "0245689abcdefgh" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("C27").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "2345689abefh" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("D27").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "0123456789adfghi" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"

$ws.Range("B28").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "0268abcdefgh" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("C28").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "2345689" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("D28").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "013456789abdghi" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"

$ws.Range("A29").HorizontalAlignment = -4108
$ws.Range("B29").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "025689abcdefgh" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("C29").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "0235689bcdeg" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"
$ws.Range("D29").Formula = "=IF(A25=" + [char]34 + [char]34 + ", " + [char]34 + [char]34 + ", IF(ISNUMBER(FIND(A25, " + [char]34 + "0123456789abceghi" + [char]34 + ")), 1, " + [char]34 + [char]34 + "))"

# --- New conditional formatting for H27 (copy of the colour-scale pair used on B18) ---
$ws.Range("H27").FormatConditions.AddColorScale(2)
$cs1 = $ws.Range("H27").FormatConditions.Item($ws.Range("H27").FormatConditions.Count)
$cs1.ColorScaleCriteria.Item(1).Type = 0
$cs1.ColorScaleCriteria.Item(1).Value = 0
$cs1.ColorScaleCriteria.Item(1).FormatColor.Color = 0xFFFFFF
$cs1.ColorScaleCriteria.Item(2).Type = 0
$cs1.ColorScaleCriteria.Item(2).Value = 1
$cs1.ColorScaleCriteria.Item(2).FormatColor.Color = 0x50D092

$ws.Range("H27").FormatConditions.AddColorScale(2)
$cs2 = $ws.Range("H27").FormatConditions.Item($ws.Range("H27").FormatConditions.Count)
$cs2.ColorScaleCriteria.Item(1).Type = -4135
$cs2.ColorScaleCriteria.Item(1).FormatColor.Color = 0xFFFCFC
$cs2.ColorScaleCriteria.Item(2).Type = -4136
$cs2.ColorScaleCriteria.Item(2).FormatColor.Color = 0x7BBE63

# --- New conditional formatting for B25:D29 (simple two-colour scale black/white) ---
$ws.Range("B25:D29").FormatConditions.AddColorScale(2)
$cs3 = $ws.Range("B25:D29").FormatConditions.Item($ws.Range("B25:D29").FormatConditions.Count)
$cs3.ColorScaleCriteria.Item(1).Type = 0
$cs3.ColorScaleCriteria.Item(1).Value = 0
$cs3.ColorScaleCriteria.Item(1).FormatColor.ThemeColor = 1
$cs3.ColorScaleCriteria.Item(2).Type = 0
$cs3.ColorScaleCriteria.Item(2).Value = 1
$cs3.ColorScaleCriteria.Item(2).FormatColor.ThemeColor = 2

# --- Selection moved to K9 as last action (reflects where the user finished editing) ---
$ws.Range("K9").Select()
